$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.994371771812439
$ws.Range("B1").Value = 2.137722253799438
$ws.Range("C1").Value = 2.147311210632324
$ws.Range("D1").Value = 2.673890352249146
$ws.Range("E1").Value = 3.431544303894043
